$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text change: "Budgets" -> "Project Name"
$ws.Range("A1").Value = "Project Name"

# Remove bold formatting from the header row
$ws.Range("A1:D1").Font.Bold = $false

# Update the active selection
[void]$ws.Range("C18").Select()
